$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new data point (row 5): date + value
$ws.Range("A5").Value = 43506
$ws.Range("B5").Value = 78055

# Move the active selection to B6 (was G5)
$ws.Range("B6").Select()
